$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.333.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6224"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07340"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2875"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07722"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.838.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001052"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6599"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.224"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.286.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "236.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.36%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.200"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9989"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.397"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1329"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06848"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.40%  "
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.479"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.013"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.935"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.152"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6786"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01818"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.783"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.232.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.657"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9427"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.991.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000118"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.876"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.809"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1126"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3852"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.39%  "
